# Fix Training Data Issue (#48)
# The "Date" column (BF) was populated with the mis-formatted string
# "6-26-2012-13" (day-month jammed against the season label) instead of
# the actual game date. Correct it to ISO format "2013-06-26" for every
# data row, keeping the value as literal text (not an auto-converted
# date serial).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldDate = "6-26-2012-13"
$newDate = "2013-06-26"

# Data rows are BF2:BF31 (row 1 is the "Date" header).
$dateRange = $ws.Range("BF2:BF31")

# Force text formatting first so Excel doesn't reinterpret the
# "yyyy-mm-dd"-shaped string as a date serial number.
$dateRange.NumberFormat = "@"
$dateRange.Value = $newDate

Write-Host "Updated Date column (BF2:BF31) from '$oldDate' to '$newDate'"
